$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sector names and recalculated average correlations for rows 4-58
$names = @('Construction & Engineering(21)', 'Air Freight & Logistics(11)', 'Marine(15)', 'Construction Materials(8)', 'Equity Real Estate Investment Trusts ...(98)', 'Trading Companies & Distributors(25)', 'Containers & Packaging(12)', 'Machinery(86)', 'Auto Components(21)', 'Multi-Utilities(18)', 'Building Products(24)', 'Chemicals(52)', 'Mortgage Real Estate Investment Trust...(16)', 'Energy Equipment & Services(38)', 'Life Sciences Tools & Services(19)', 'Wireless Telecommunication Services(14)', 'Airlines(14)', 'Insurance(75)', 'Gas Utilities(12)', 'IT Services(52)', 'Leisure Products(11)', 'Semiconductors & Semiconductor Equipment(70)', 'Electrical Equipment(28)', 'Metals & Mining(106)', 'Capital Markets(76)', 'Household Durables(39)', 'Oil, Gas & Consumable Fuels(125)', 'Technology Hardware, Storage & Periph...(19)', 'Water Utilities(13)', 'Professional Services(35)', 'Electric Utilities(28)', 'ETF(303)', 'Health Care Providers & Services(47)', 'Communications Equipment(45)', 'Banks(251)', 'Consumer Finance(15)', 'Aerospace & Defense(37)', 'Specialty Retail(59)', 'Hotels, Restaurants & Leisure(51)', 'Electronic Equipment, Instruments & C...(78)', 'Commercial Services & Supplies(52)', 'Textiles, Apparel & Luxury Goods(29)', 'Software(70)', 'Beverages(21)', 'Diversified Consumer Services(17)', 'Food & Staples Retailing(16)', 'Real Estate Management & Development(23)', 'Entertainment(22)', 'Media(42)', 'Diversified Telecommunication Services(20)', 'Health Care Equipment & Supplies(86)', 'Food Products(46)', 'Thrifts & Mortgage Finance(47)', 'Biotechnology(128)', 'Pharmaceuticals(53)')
$vals = @(0.5786025299199333, 0.5775675828705409, 0.5690460143460049, 0.5663544204918297, 0.5496987913628164, 0.5441324810134711, 0.5275944728785374, 0.5235453611083485, 0.5046759631775736, 0.4909455662030818, 0.4865975412683555, 0.4779022354291935, 0.4726946262294544, 0.4720797566113515, 0.4662222783206468, 0.4456355275386176, 0.4420316194123448, 0.4331316156024389, 0.4197610392852889, 0.4142540919728485, 0.4131002866899252, 0.4113315235717009, 0.4110709898800927, 0.4107914148021186, 0.404324614074526, 0.4036972876281053, 0.4020471117340391, 0.3981848729953394, 0.380142963547861, 0.3793475327116087, 0.3626800602595638, 0.3522363555472797, 0.3459392702200678, 0.3451938496260053, 0.3362675392680736, 0.3360127359773946, 0.3327326980024479, 0.3314211153553554, 0.330966902539903, 0.3233749530919108, 0.3159174398305094, 0.3084588888265488, 0.3025465012445216, 0.3010316192346816, 0.2979870840024118, 0.2960646697541379, 0.2727293729642281, 0.2468887503039209, 0.2466866844775422, 0.2287517233842288, 0.2240240279225028, 0.2119574383509444, 0.1785330825979096, 0.1757997545802315, 0.1446982411768738)

$startRow = 4
for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 2).Value = $vals[$i]
}

